$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) for columns L, M, N
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# New data values for rows 2-7
$values = @(
    @(92.42558549437925, 211586, 306.6463768115942),
    @(86.73427165141923, 51967, 371.1928571428571),
    @(88.21289274334562, 150490, 137.0582877959927),
    @(95.69895606541161, 53156, 186.5122807017544),
    @(19.42841748106926, 2114, 14.99290780141844),
    @(22.50301276387648, 94, 8.545454545454545)
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Range("L$row").Value = $values[$i][0]
    $ws.Range("M$row").Value = $values[$i][1]
    $ws.Range("N$row").Value = $values[$i][2]
}
